$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: set a cell to an exact text value, even if it looks numeric,
# matching the source data which stores Price/Volume figures as text.
function Set-TextValue($cell, $text) {
    $cell.NumberFormat = "@"
    $cell.Value = $text
    $cell.Style = "Normal"
}

# Rows 2-48, 51: refreshed Price (D) and Volume(1h) (E) figures
$ws.Range("D2").Value = "67.548.95"
$ws.Range("E2").Value = "  +0.16%  "
$ws.Range("D3").Value = "3.507.10"
$ws.Range("E3").Value = "  -0.31%  "
Set-TextValue $ws.Range("D4") "1.00"
$ws.Range("E4").Value = "  +0.04%  "
Set-TextValue $ws.Range("D5") "609.79"
$ws.Range("E5").Value = "  -0.16%  "
Set-TextValue $ws.Range("D6") "152.30"
$ws.Range("E6").Value = "  +1.15%  "
$ws.Range("D7").Value = "3.504.62"
$ws.Range("E7").Value = "  -0.34%  "
$ws.Range("E8").Value = "  -0.05%  "
$ws.Range("E9").Value = "  +1.27%  "
$ws.Range("E10").Value = "  +3.01%  "
$ws.Range("E11").Value = "  +8.56%  "
$ws.Range("E12").Value = "  +1.63%  "
Set-TextValue $ws.Range("D13") "32.67"
$ws.Range("E13").Value = "  +2.74%  "
$ws.Range("E14").Value = "  -1.37%  "
$ws.Range("D15").Value = "4.099.04"
$ws.Range("E15").Value = "  -0.37%  "
$ws.Range("D16").Value = "3.505.36"
$ws.Range("E16").Value = "  -0.37%  "
$ws.Range("D17").Value = "67.448.50"
$ws.Range("E17").Value = "  +0.03%  "
$ws.Range("E18").Value = "  +0.03%  "
$ws.Range("E19").Value = "  +2.88%  "
Set-TextValue $ws.Range("D20") "15.57"
$ws.Range("E20").Value = "  +2.22%  "
Set-TextValue $ws.Range("D21") "9.91"
$ws.Range("E21").Value = "  +7.03%  "
Set-TextValue $ws.Range("D22") "448.18"
$ws.Range("E22").Value = "  +1.11%  "
Set-TextValue $ws.Range("D23") "0.633"
$ws.Range("E23").Value = "  +1.36%  "
Set-TextValue $ws.Range("D24") "78.18"
$ws.Range("E24").Value = "  +1.20%  "
$ws.Range("D25").Value = "3.646.35"
$ws.Range("E25").Value = "  -0.35%  "
$ws.Range("E26").Value = "  -0.01%  "
Set-TextValue $ws.Range("D27") "0.0000126"
$ws.Range("E27").Value = "  -2.40%  "
Set-TextValue $ws.Range("D28") "8.85"
$ws.Range("E28").Value = "  +6.01%  "
Set-TextValue $ws.Range("D29") "10.11"
$ws.Range("E29").Value = "  -0.88%  "
$ws.Range("E30").Value = "  +0.79%  "
Set-TextValue $ws.Range("D31") "1.67"
$ws.Range("E31").Value = "  +7.93%  "
Set-TextValue $ws.Range("D32") "0.170"
$ws.Range("E32").Value = "  +3.71%  "
$ws.Range("E33").Value = "  -0.01%  "
Set-TextValue $ws.Range("D34") "25.73"
$ws.Range("E34").Value = "  -0.20%  "
Set-TextValue $ws.Range("D35") "6.19"
$ws.Range("E35").Value = "  +0.90%  "
$ws.Range("E36").Value = "  +1.81%  "
$ws.Range("D37").Value = "3.498.12"
$ws.Range("E37").Value = "  -0.43%  "
Set-TextValue $ws.Range("D38") "8.06"
$ws.Range("E38").Value = "  +0.38%  "
$ws.Range("E40").Value = "  +6.44%  "
Set-TextValue $ws.Range("D41") "1.00"
$ws.Range("E41").Value = "  +0.00%  "
Set-TextValue $ws.Range("D42") "0.0897"
$ws.Range("E42").Value = "  +2.70%  "
Set-TextValue $ws.Range("D43") "173.60"
$ws.Range("E43").Value = "  -2.58%  "
Set-TextValue $ws.Range("D44") "5.50"
$ws.Range("E44").Value = "  +1.23%  "
Set-TextValue $ws.Range("D45") "30.11"
$ws.Range("E45").Value = "  +8.61%  "
$ws.Range("E46").Value = "  +0.37%  "
Set-TextValue $ws.Range("D47") "46.83"
$ws.Range("E47").Value = "  +2.83%  "
$ws.Range("E48").Value = "  +3.67%  "
Set-TextValue $ws.Range("D51") "0.254"
$ws.Range("E51").Value = "  +3.36%  "

# Rows 49-50: dogwifhat and Cosmos swapped position, with refreshed Price/Volume figures
$ws.Range("B49").Value = "dogwifhat"
$ws.Range("C49").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
Set-TextValue $ws.Range("D49") "2.53"
$ws.Range("E49").Value = "  -2.52%  "
$ws.Range("B50").Value = "Cosmos"
$ws.Range("C50").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
Set-TextValue $ws.Range("D50") "7.67"
$ws.Range("E50").Value = "  +1.33%  "
